$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2499.5
$ws.Range("J48").Value = 2499.5
$ws.Range("L48").Value = 7498.5
$ws.Range("N48").Value = -8082.5
$ws.Range("H56").Value = 2499.5
$ws.Range("J56").Value = 2499.5
$ws.Range("L56").Value = 7498.5
$ws.Range("N56").Value = -8566.5
$ws.Range("H64").Value = 3518.2727
$ws.Range("I64").Value = 3530.5789
$ws.Range("J64").Value = 3508.92
$ws.Range("K64").Value = 3530.5789
$ws.Range("L64").Value = 3508.92
$ws.Range("M64").Value = -3282.5789
$ws.Range("N64").Value = -4004.92
$ws.Range("H67").Value = 3518.2727
$ws.Range("I67").Value = 3530.5789
$ws.Range("J67").Value = 3508.92
$ws.Range("K67").Value = 3530.5789
$ws.Range("L67").Value = 3508.92
$ws.Range("M67").Value = -2672.5789
$ws.Range("N67").Value = -5224.92
$ws.Range("H116").Value = 3425.2942
$ws.Range("I116").Value = 2656.9092
$ws.Range("K116").Value = 2656.9092
$ws.Range("M116").Value = 785.0907999999999
$ws.Range("H135").Value = 249.42857
$ws.Range("I135").Value = 226
$ws.Range("K135").Value = 2034
$ws.Range("M135").Value = 501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2667.7666
$ws.Range("I32").Value = 2788.4717
$ws.Range("K32").Value = 2788.4717
$ws.Range("M32").Value = -2501.4717
$ws.Range("H74").Value = 981.5714
$ws.Range("J74").Value = 547
$ws.Range("L74").Value = 547
$ws.Range("N74").Value = -2295
$ws.Range("H77").Value = 981.5714
$ws.Range("J77").Value = 547
$ws.Range("L77").Value = 2735
$ws.Range("N77").Value = -11471
$ws.Range("H110").Value = 1316.6897
$ws.Range("I110").Value = 1176.4348
$ws.Range("J110").Value = 1854.3334
$ws.Range("K110").Value = 1176.4348
$ws.Range("L110").Value = 1854.3334
$ws.Range("M110").Value = 868.5652
$ws.Range("N110").Value = -5944.3334
$ws.Range("H122").Value = 1569
$ws.Range("I122").Value = 2100
$ws.Range("J122").Value = 1038
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 3114
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -8014
$ws.Range("H132").Value = 3189.5
$ws.Range("I132").Value = 2939.1538
$ws.Range("K132").Value = 8817.4614
$ws.Range("M132").Value = -6287.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 333335500
$ws.Range("I105").Value = 500002500
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 500002500
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -500000753
$ws.Range("N105").Value = -4994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N7").Value = ""
$ws.Range("H7").Value = 397
$ws.Range("I7").Value = 397
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 397
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -284
$ws.Range("H22").Value = 469.0909
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50
$ws.Range("H31").Value = 1221.3103
$ws.Range("I31").Value = 851.9211
$ws.Range("K31").Value = 851.9211
$ws.Range("M31").Value = -556.9211
$ws.Range("H34").Value = 1221.3103
$ws.Range("I34").Value = 851.9211
$ws.Range("K34").Value = 851.9211
$ws.Range("M34").Value = -649.9211
$ws.Range("H55").Value = 7300
$ws.Range("I55").Value = 7300
$ws.Range("K55").Value = 7300
$ws.Range("M55").Value = -6985
$ws.Range("H62").Value = 5408341
$ws.Range("I62").Value = 2981.2188
$ws.Range("J62").Value = 40002644
$ws.Range("K62").Value = 2981.2188
$ws.Range("L62").Value = 40002644
$ws.Range("M62").Value = -2357.2188
$ws.Range("N62").Value = -40003892
$ws.Range("H65").Value = 5408341
$ws.Range("I65").Value = 2981.2188
$ws.Range("J65").Value = 40002644
$ws.Range("K65").Value = 14906.094
$ws.Range("L65").Value = 200013220
$ws.Range("M65").Value = -11786.094
$ws.Range("N65").Value = -200019460
$ws.Range("N122").Value = ""
$ws.Range("H122").Value = 826.8
$ws.Range("I122").Value = 826.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2480.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -30.39999999999964
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("M131").Value = ""
$ws.Range("N131").Value = ""
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 11231.9
$ws.Range("I132").Value = 13790
$ws.Range("J132").Value = 999.5
$ws.Range("K132").Value = 41370
$ws.Range("L132").Value = 2998.5
$ws.Range("M132").Value = -38840
$ws.Range("N132").Value = -8058.5
$ws.Range("H134").Value = 9260736
$ws.Range("I134").Value = 13334787
$ws.Range("J134").Value = 1528
$ws.Range("K134").Value = 40004361
$ws.Range("L134").Value = 4584
$ws.Range("M134").Value = -40001826
$ws.Range("N134").Value = -9654

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2200
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 2300
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 6900
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -7466
$ws.Range("H68").Value = 2164.62
$ws.Range("J68").Value = 2198.5918
$ws.Range("L68").Value = 6595.7754
$ws.Range("N68").Value = -8217.7754
$ws.Range("H71").Value = 2164.62
$ws.Range("J71").Value = 2198.5918
$ws.Range("L71").Value = 19787.3262
$ws.Range("N71").Value = -27899.3262
$ws.Range("H107").Value = 7981.7856
$ws.Range("I107").Value = 1426.5
$ws.Range("J107").Value = 9074.333000000001
$ws.Range("K107").Value = 4279.5
$ws.Range("L107").Value = 27222.999
$ws.Range("M107").Value = -2359.5
$ws.Range("N107").Value = -31062.999
$ws.Range("H122").Value = 581.2
$ws.Range("I122").Value = 427.375
$ws.Range("K122").Value = 3846.375
$ws.Range("M122").Value = -1396.375
$ws.Range("H131").Value = 43480160
$ws.Range("J131").Value = 2232
$ws.Range("L131").Value = 6696
$ws.Range("N131").Value = -16776
$ws.Range("H137").Value = 14017.091
$ws.Range("J137").Value = 23516.334
$ws.Range("L137").Value = 70549.00199999999
$ws.Range("N137").Value = -80749.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 24997.5
$ws.Range("J26").Value = 24997.5
$ws.Range("L26").Value = 24997.5
$ws.Range("N26").Value = -25557.5
$ws.Range("H50").Value = 24997.5
$ws.Range("J50").Value = 24997.5
$ws.Range("L50").Value = 24997.5
$ws.Range("N50").Value = -25993.5
$ws.Range("H80").Value = 4592
$ws.Range("I80").Value = 3338.125
$ws.Range("J80").Value = 7099.75
$ws.Range("K80").Value = 3338.125
$ws.Range("L80").Value = 7099.75
$ws.Range("M80").Value = -2340.125
$ws.Range("N80").Value = -9095.75
$ws.Range("H83").Value = 4592
$ws.Range("I83").Value = 3338.125
$ws.Range("J83").Value = 7099.75
$ws.Range("K83").Value = 16690.625
$ws.Range("L83").Value = 35498.75
$ws.Range("M83").Value = -11698.625
$ws.Range("N83").Value = -45482.75
$ws.Range("H107").Value = 916438.2
$ws.Range("I107").Value = 1748890.9
$ws.Range("J107").Value = 740.3
$ws.Range("K107").Value = 1748890.9
$ws.Range("L107").Value = 740.3
$ws.Range("M107").Value = -1746970.9
$ws.Range("N107").Value = -4580.3
$ws.Range("H132").Value = 3151.682
$ws.Range("I132").Value = 2790
$ws.Range("J132").Value = 4116.1665
$ws.Range("K132").Value = 8370
$ws.Range("L132").Value = 12348.4995
$ws.Range("M132").Value = -5840
$ws.Range("N132").Value = -17408.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M45").Value = ""
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("H55").Value = 516.7273
$ws.Range("I55").Value = 89.333336
$ws.Range("J55").Value = 1029.6
$ws.Range("K55").Value = 89.333336
$ws.Range("L55").Value = 1029.6
$ws.Range("M55").Value = 83.666664
$ws.Range("N55").Value = -1375.6
$ws.Range("H132").Value = 29742.5
$ws.Range("J132").Value = 57822.445
$ws.Range("L132").Value = 173467.335
$ws.Range("N132").Value = -178527.335
$ws.Range("H136").Value = 10979.272
$ws.Range("I136").Value = 15642.286
$ws.Range("K136").Value = 46926.858
$ws.Range("M136").Value = -44376.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26323040
$ws.Range("I62").Value = 35719588
$ws.Range("K62").Value = 35719588
$ws.Range("M62").Value = -35718964
$ws.Range("H65").Value = 26323040
$ws.Range("I65").Value = 35719588
$ws.Range("K65").Value = 178597940
$ws.Range("M65").Value = -178594820
$ws.Range("H132").Value = 5388.3794
$ws.Range("I132").Value = 6315.1113
$ws.Range("J132").Value = 3871.9092
$ws.Range("K132").Value = 18945.3339
$ws.Range("L132").Value = 11615.7276
$ws.Range("M132").Value = -16415.3339
$ws.Range("N132").Value = -16675.7276
$ws.Range("H138").Value = 31643
$ws.Range("J138").Value = 31643
$ws.Range("L138").Value = 31643
$ws.Range("N138").Value = -41923
